$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "チーム" (Team) -> "コメント" (Comment)
$ws.Range("B2").Value = "コメント"

# Comment column (B) per staff row: replace old "A"/"B" team labels with
# real free-text comments, or clear the cell where no comment applies.
$ws.Range("B5").Value = "夜勤4回まで"
$ws.Range("B6").Value = "夜勤4回まで"
$ws.Range("B7").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("B9").ClearContents()
$ws.Range("B10").ClearContents()
$ws.Range("B11").ClearContents()
$ws.Range("B12").Value = "夜勤3回まで"
$ws.Range("B13").ClearContents()
$ws.Range("B14").ClearContents()
$ws.Range("B15").ClearContents()
$ws.Range("B16").ClearContents()
$ws.Range("B17").ClearContents()
$ws.Range("B18").ClearContents()
$ws.Range("B19").Value = "新人　月前半長夜勤なし　"
$ws.Range("B20").ClearContents()
$ws.Range("B21").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("B23").Value = "土日休み日勤のみ"
$ws.Range("B24").Value = "夜勤土日のみ3回まで"
$ws.Range("B25").Value = "長入明　水木金3回まで"
$ws.Range("B26").ClearContents()
$ws.Range("B27").ClearContents()
$ws.Range("B28").ClearContents()
$ws.Range("B29").ClearContents()
$ws.Range("B30").ClearContents()
$ws.Range("B31").ClearContents()
$ws.Range("B32").ClearContents()
$ws.Range("B33").Value = "新人　月前半長夜勤なし"

# Columns C:G, rows 4-33: restyle from the old bordered/white cell style to
# match the plain beige style already used in columns H onward (copy format
# only from an already-correctly-styled cell so the cell *values* are left
# untouched).
$ws.Range("H5").Copy()
$ws.Range("C4:G33").PasteSpecial(-4122)
